$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "skeleton_archer_blue"
$ws.Range("B54").Value = 52
$ws.Range("C54").Value = "Prefabs/Object/NPC/skeleton_archer_blue"
$ws.Range("D54").Value = 2
$ws.Range("E54").Value = 20
$ws.Range("F54").Value = "DropBag_1"
$ws.Range("H54").Value = 2
$ws.Range("I54").NumberFormat = "@"
$ws.Range("I54").Value = "PlayerAtt52"
$ws.Range("J54").Value = "ConsumeData_1"

$ws.Range("A55").Value = "skeleton_archer_green"
$ws.Range("C55").Value = "Prefabs/Object/NPC/skeleton_archer_green"
$ws.Range("D55").Value = 2
$ws.Range("E55").Value = 20
$ws.Range("F55").Value = "DropBag_1"
$ws.Range("H55").Value = 2
$ws.Range("I55").NumberFormat = "@"
$ws.Range("I55").Value = "PlayerAtt52"
$ws.Range("J55").Value = "ConsumeData_1"

$ws.Range("A56").Value = "skeleton_archer_purple"
$ws.Range("C56").Value = "Prefabs/Object/NPC/skeleton_archer_purple"
$ws.Range("D56").Value = 2
$ws.Range("E56").Value = 20
$ws.Range("F56").Value = "DropBag_1"
$ws.Range("H56").Value = 2
$ws.Range("I56").NumberFormat = "@"
$ws.Range("I56").Value = "PlayerAtt52"
$ws.Range("J56").Value = "ConsumeData_1"

$ws.Range("A57").Value = "skeleton_archer_red"
$ws.Range("C57").Value = "Prefabs/Object/NPC/skeleton_archer_red"
$ws.Range("D57").Value = 2
$ws.Range("E57").Value = 20
$ws.Range("F57").Value = "DropBag_1"
$ws.Range("H57").Value = 2
$ws.Range("I57").NumberFormat = "@"
$ws.Range("I57").Value = "PlayerAtt52"
$ws.Range("J57").Value = "ConsumeData_1"

$ws.Range("A58").Value = "skeleton_archer_teal"
$ws.Range("C58").Value = "Prefabs/Object/NPC/skeleton_archer_teal"
$ws.Range("D58").Value = 2
$ws.Range("E58").Value = 20
$ws.Range("F58").Value = "DropBag_1"
$ws.Range("H58").Value = 2
$ws.Range("I58").NumberFormat = "@"
$ws.Range("I58").Value = "PlayerAtt52"
$ws.Range("J58").Value = "ConsumeData_1"

$ws.Range("A59").Value = "skeleton_archer_yellow"
$ws.Range("C59").Value = "Prefabs/Object/NPC/skeleton_archer_yellow"
$ws.Range("D59").Value = 2
$ws.Range("E59").Value = 20
$ws.Range("F59").Value = "DropBag_1"
$ws.Range("H59").Value = 2
$ws.Range("I59").NumberFormat = "@"
$ws.Range("I59").Value = "PlayerAtt52"
$ws.Range("J59").Value = "ConsumeData_1"

$ws.Range("A60").Value = "skeleton_tom_angry"
$ws.Range("C60").Value = "Prefabs/Object/NPC/skeleton_tom_angry"
$ws.Range("D60").Value = 2
$ws.Range("E60").Value = 20
$ws.Range("F60").Value = "DropBag_1"
$ws.Range("H60").Value = 2
$ws.Range("I60").NumberFormat = "@"
$ws.Range("I60").Value = "PlayerAtt52"
$ws.Range("J60").Value = "ConsumeData_1"

$ws.Range("A61").Value = "skeleton_tom_happy"
$ws.Range("C61").Value = "Prefabs/Object/NPC/skeleton_tom_happy"
$ws.Range("D61").Value = 2
$ws.Range("E61").Value = 20
$ws.Range("F61").Value = "DropBag_1"
$ws.Range("H61").Value = 2
$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = "PlayerAtt52"
$ws.Range("J61").Value = "ConsumeData_1"

$ws.Range("A62").Value = "skeleton_king_blue"
$ws.Range("C62").Value = "Prefabs/Object/NPC/skeleton_king_blue"
$ws.Range("D62").Value = 2
$ws.Range("E62").Value = 20
$ws.Range("F62").Value = "DropBag_1"
$ws.Range("H62").Value = 2
$ws.Range("I62").NumberFormat = "@"
$ws.Range("I62").Value = "PlayerAtt52"
$ws.Range("J62").Value = "ConsumeData_1"

$ws.Range("A63").Value = "skeleton_king_green"
$ws.Range("C63").Value = "Prefabs/Object/NPC/skeleton_king_green"
$ws.Range("D63").Value = 2
$ws.Range("E63").Value = 20
$ws.Range("F63").Value = "DropBag_1"
$ws.Range("H63").Value = 2
$ws.Range("I63").NumberFormat = "@"
$ws.Range("I63").Value = "PlayerAtt52"
$ws.Range("J63").Value = "ConsumeData_1"

$ws.Range("A64").Value = "skeleton_king_purple"
$ws.Range("C64").Value = "Prefabs/Object/NPC/skeleton_king_purple"
$ws.Range("D64").Value = 2
$ws.Range("E64").Value = 20
$ws.Range("F64").Value = "DropBag_1"
$ws.Range("H64").Value = 2
$ws.Range("I64").NumberFormat = "@"
$ws.Range("I64").Value = "PlayerAtt52"
$ws.Range("J64").Value = "ConsumeData_1"

$ws.Range("A65").Value = "skeleton_king_red"
$ws.Range("C65").Value = "Prefabs/Object/NPC/skeleton_king_red"
$ws.Range("D65").Value = 2
$ws.Range("E65").Value = 20
$ws.Range("F65").Value = "DropBag_1"
$ws.Range("H65").Value = 2
$ws.Range("I65").NumberFormat = "@"
$ws.Range("I65").Value = "PlayerAtt52"
$ws.Range("J65").Value = "ConsumeData_1"

$ws.Range("A66").Value = "skeleton_king_teal"
$ws.Range("C66").Value = "Prefabs/Object/NPC/skeleton_king_teal"
$ws.Range("D66").Value = 2
$ws.Range("E66").Value = 20
$ws.Range("F66").Value = "DropBag_1"
$ws.Range("H66").Value = 2
$ws.Range("I66").NumberFormat = "@"
$ws.Range("I66").Value = "PlayerAtt52"
$ws.Range("J66").Value = "ConsumeData_1"

$ws.Range("A67").Value = "skeleton_king_yellow"
$ws.Range("C67").Value = "Prefabs/Object/NPC/skeleton_king_yellow"
$ws.Range("D67").Value = 2
$ws.Range("E67").Value = 20
$ws.Range("F67").Value = "DropBag_1"
$ws.Range("H67").Value = 2
$ws.Range("I67").NumberFormat = "@"
$ws.Range("I67").Value = "PlayerAtt52"
$ws.Range("J67").Value = "ConsumeData_1"

$ws.Range("A68").Value = "skeleton_mage_blue"
$ws.Range("C68").Value = "Prefabs/Object/NPC/skeleton_mage_blue"
$ws.Range("D68").Value = 2
$ws.Range("E68").Value = 20
$ws.Range("F68").Value = "DropBag_1"
$ws.Range("H68").Value = 2
$ws.Range("I68").NumberFormat = "@"
$ws.Range("I68").Value = "PlayerAtt52"
$ws.Range("J68").Value = "ConsumeData_1"

$ws.Range("A69").Value = "skeleton_mage_green"
$ws.Range("C69").Value = "Prefabs/Object/NPC/skeleton_mage_green"
$ws.Range("D69").Value = 2
$ws.Range("E69").Value = 20
$ws.Range("F69").Value = "DropBag_1"
$ws.Range("H69").Value = 2
$ws.Range("I69").NumberFormat = "@"
$ws.Range("I69").Value = "PlayerAtt52"
$ws.Range("J69").Value = "ConsumeData_1"

$ws.Range("A70").Value = "skeleton_mage_purple"
$ws.Range("C70").Value = "Prefabs/Object/NPC/skeleton_mage_purple"
$ws.Range("D70").Value = 2
$ws.Range("E70").Value = 20
$ws.Range("F70").Value = "DropBag_1"
$ws.Range("H70").Value = 2
$ws.Range("I70").NumberFormat = "@"
$ws.Range("I70").Value = "PlayerAtt52"
$ws.Range("J70").Value = "ConsumeData_1"

$ws.Range("A71").Value = "skeleton_mage_red"
$ws.Range("C71").Value = "Prefabs/Object/NPC/skeleton_mage_red"
$ws.Range("D71").Value = 2
$ws.Range("E71").Value = 20
$ws.Range("F71").Value = "DropBag_1"
$ws.Range("H71").Value = 2
$ws.Range("I71").NumberFormat = "@"
$ws.Range("I71").Value = "PlayerAtt52"
$ws.Range("J71").Value = "ConsumeData_1"

$ws.Range("A72").Value = "skeleton_mage_teal"
$ws.Range("C72").Value = "Prefabs/Object/NPC/skeleton_mage_teal"
$ws.Range("D72").Value = 2
$ws.Range("E72").Value = 20
$ws.Range("F72").Value = "DropBag_1"
$ws.Range("H72").Value = 2
$ws.Range("I72").NumberFormat = "@"
$ws.Range("I72").Value = "PlayerAtt52"
$ws.Range("J72").Value = "ConsumeData_1"

$ws.Range("A73").Value = "skeleton_mage_yellow"
$ws.Range("C73").Value = "Prefabs/Object/NPC/skeleton_mage_yellow"
$ws.Range("D73").Value = 2
$ws.Range("E73").Value = 20
$ws.Range("F73").Value = "DropBag_1"
$ws.Range("H73").Value = 2
$ws.Range("I73").NumberFormat = "@"
$ws.Range("I73").Value = "PlayerAtt52"
$ws.Range("J73").Value = "ConsumeData_1"

$ws.Range("A74").Value = "skeleton_warrior_blue"
$ws.Range("C74").Value = "Prefabs/Object/NPC/skeleton_warrior_blue"
$ws.Range("D74").Value = 2
$ws.Range("E74").Value = 20
$ws.Range("F74").Value = "DropBag_1"
$ws.Range("H74").Value = 2
$ws.Range("I74").NumberFormat = "@"
$ws.Range("I74").Value = "PlayerAtt52"
$ws.Range("J74").Value = "ConsumeData_1"

$ws.Range("A75").Value = "skeleton_warrior_green"
$ws.Range("C75").Value = "Prefabs/Object/NPC/skeleton_warrior_green"
$ws.Range("D75").Value = 2
$ws.Range("E75").Value = 20
$ws.Range("F75").Value = "DropBag_1"
$ws.Range("H75").Value = 2
$ws.Range("I75").NumberFormat = "@"
$ws.Range("I75").Value = "PlayerAtt52"
$ws.Range("J75").Value = "ConsumeData_1"

$ws.Range("A76").Value = "skeleton_warrior_purple"
$ws.Range("C76").Value = "Prefabs/Object/NPC/skeleton_warrior_purple"
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 20
$ws.Range("F76").Value = "DropBag_1"
$ws.Range("H76").Value = 2
$ws.Range("I76").NumberFormat = "@"
$ws.Range("I76").Value = "PlayerAtt52"
$ws.Range("J76").Value = "ConsumeData_1"

$ws.Range("A77").Value = "skeleton_warrior_red"
$ws.Range("C77").Value = "Prefabs/Object/NPC/skeleton_warrior_red"
$ws.Range("D77").Value = 2
$ws.Range("E77").Value = 20
$ws.Range("F77").Value = "DropBag_1"
$ws.Range("H77").Value = 2
$ws.Range("I77").NumberFormat = "@"
$ws.Range("I77").Value = "PlayerAtt52"
$ws.Range("J77").Value = "ConsumeData_1"

$ws.Range("A78").Value = "skeleton_warrior_teal"
$ws.Range("C78").Value = "Prefabs/Object/NPC/skeleton_warrior_teal"
$ws.Range("D78").Value = 2
$ws.Range("E78").Value = 20
$ws.Range("F78").Value = "DropBag_1"
$ws.Range("H78").Value = 2
$ws.Range("I78").NumberFormat = "@"
$ws.Range("I78").Value = "PlayerAtt52"
$ws.Range("J78").Value = "ConsumeData_1"

$ws.Range("A79").Value = "skeleton_warrior_yellow"
$ws.Range("C79").Value = "Prefabs/Object/NPC/skeleton_warrior_yellow"
$ws.Range("D79").Value = 2
$ws.Range("E79").Value = 20
$ws.Range("F79").Value = "DropBag_1"
$ws.Range("H79").Value = 2
$ws.Range("I79").NumberFormat = "@"
$ws.Range("I79").Value = "PlayerAtt52"
$ws.Range("J79").Value = "ConsumeData_1"

$ws.Range("J53:J79").Select()
